$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report for "Rabanito" at Vega Modelo de Temuco was
# inserted as the new top data row (row 4, right after the two anchor
# rows). This pushes all existing data rows (old 4-108) down by one,
# which is why every subsequent row's values shift down by one position.
$ws.Rows("4").Insert()

$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 45083
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 300000001
$ws.Range("G4").Value = "Rabanito"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7667
$ws.Range("N4").Value = "$/docena de paquetes"
$ws.Range("O4").Value = "Provincia de Cautín"
$ws.Range("P4").Value = 639
$ws.Range("Q4").Value = 12
$ws.Range("R4").Value = "Hortaliza"
